$wb = $excel.ActiveWorkbook

# 展览 (sheet1)
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 1175
$ws.Cells.Item(3, 7).Value = 65
$ws.Cells.Item(4, 6).Value = 1608
$ws.Cells.Item(4, 7).Value = 70
$ws.Cells.Item(5, 6).Value = 183
$ws.Cells.Item(6, 6).Value = 183
$ws.Cells.Item(7, 6).Value = 38
$ws.Cells.Item(8, 6).Value = 1563
$ws.Cells.Item(9, 6).Value = 3170
$ws.Cells.Item(10, 6).Value = 711
$ws.Cells.Item(11, 6).Value = 1897
$ws.Cells.Item(12, 6).Value = 1853
$ws.Cells.Item(13, 6).Value = 919
$ws.Cells.Item(15, 6).Value = 10
$ws.Cells.Item(16, 6).Value = 1537
$ws.Cells.Item(17, 6).Value = 311
$ws.Cells.Item(19, 6).Value = 43
$ws.Cells.Item(20, 6).Value = 1329
$ws.Cells.Item(21, 6).Value = 443
$ws.Cells.Item(22, 6).Value = 545
$ws.Cells.Item(23, 6).Value = 233
$ws.Cells.Item(24, 6).Value = 7909
$ws.Cells.Item(25, 6).Value = 9239
$ws.Cells.Item(26, 6).Value = 796
$ws.Cells.Item(27, 6).Value = 608
$ws.Cells.Item(28, 6).Value = 1755
$ws.Cells.Item(29, 6).Value = 113
$ws.Cells.Item(30, 6).Value = 290

# 演出 (sheet2)
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(5, 6).Value = 108

# 本地生活 (sheet3)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 54

# 全部类型 (sheet4)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 54
$ws.Cells.Item(4, 6).Value = 1175
$ws.Cells.Item(4, 7).Value = 65
$ws.Cells.Item(5, 6).Value = 1608
$ws.Cells.Item(5, 7).Value = 70
$ws.Cells.Item(6, 6).Value = 183
$ws.Cells.Item(7, 6).Value = 183
$ws.Cells.Item(9, 6).Value = 38
$ws.Cells.Item(10, 6).Value = 1563
$ws.Cells.Item(11, 6).Value = 3170
$ws.Cells.Item(12, 6).Value = 711
$ws.Cells.Item(13, 6).Value = 1897
$ws.Cells.Item(14, 6).Value = 1853
$ws.Cells.Item(15, 6).Value = 919
$ws.Cells.Item(17, 6).Value = 10
$ws.Cells.Item(18, 6).Value = 1537
$ws.Cells.Item(19, 6).Value = 311
$ws.Cells.Item(22, 6).Value = 43
$ws.Cells.Item(24, 6).Value = 1329
$ws.Cells.Item(25, 6).Value = 443
$ws.Cells.Item(26, 6).Value = 545
$ws.Cells.Item(27, 6).Value = 233
$ws.Cells.Item(28, 6).Value = 7909
$ws.Cells.Item(29, 6).Value = 9239
$ws.Cells.Item(30, 6).Value = 796
$ws.Cells.Item(31, 6).Value = 608
$ws.Cells.Item(32, 6).Value = 1755
$ws.Cells.Item(33, 6).Value = 108
$ws.Cells.Item(35, 6).Value = 113
$ws.Cells.Item(36, 6).Value = 290
